$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @(0.8080807042508411, 0.1644293338300642, 0.1433246603658134, 2.08772282293279, 0.00252037125172649, 0.2403379213716903, 0.3867118761887127, 0.3154163334477289, 5.470469219942089),
    @(0.7694453040819269, 0.1643014535607001, 0.141205728894505, 2.096430486766238, 0.002523015945694305, 0.2415781363556349, 0.3526390045186929, 0.3105214654573416, 5.501896778246902),
    @(0.7460048284323761, 0.1642360130567013, 0.1399557901362485, 2.102684395150035, 0.002524727336564507, 0.2424222113708971, 0.3317872330759712, 0.3076400480756192, 5.523603024747928),
    @(0.7365242307037079, 0.1642126562354846, 0.1394593493931424, 2.105461277789409, 0.002525446819790089, 0.242786972454816, 0.3233078777176814, 0.3064971765733873, 5.533054731253188),
    @(0.7349543296014645, 0.1642089784622982, 0.1393776981947994, 2.105936175536371, 0.002525567624597058, 0.242848797487099, 0.3219009847162937, 0.3063092998541919, 5.534660804493399),
    @(0.7458766789663969, 0.1642356846243729, 0.1399490425571486, 2.102720920259848, 0.002524736950479267, 0.2424270464354912, 0.3316728042959625, 0.3076245078783089, 5.52372803872521),
    @(0.7947011460056785, 0.1643825392946532, 0.142583490218847, 2.090537018643317, 0.002521265012766789, 0.2407484255960064, 0.3749495143986508, 0.3137029153087525, 5.480805519456908),
    @(0.8926556043767562, 0.1647733794872934, 0.1481522754924782, 2.073836320730962, 0.002515148161920542, 0.2381106556054426, 0.4603443520740598, 0.3266018712745762, 5.415741505964405),
    @(0.965942634188707, 0.1651220864001886, 0.1524859632556712, 2.065942019120712, 0.002511071574108734, 0.2365698088022157, 0.5233875553033727, 0.3366701830089056, 5.379573629572235),
    @(0.9995641398391797, 0.1652938537070128, 0.15450944980131, 2.063299127300439, 0.002509306792139741, 0.2359547429604589, 0.5521298434086077, 0.3413778631244355, 5.365643701282949),
    @(1.012335790324272, 0.1653607664180967, 0.1552831146413922, 2.062434527776745, 0.002508651343890116, 0.2357341550016159, 0.5630225019697548, 0.3431787630561729, 5.36073136835094),
    @(1.009583423682642, 0.1653462728161372, 0.1551161634825604, 2.06261467938441, 0.002508791936474256, 0.2357811148329318, 0.560676200611681, 0.3427900996261855, 5.361773201977229),
    @(1.000614075580586, 0.1652993213356808, 0.1545729514965188, 2.063225267278895, 0.002509252611013677, 0.2359363482016867, 0.5530258203307881, 0.3415256602108485, 5.36523229410102),
    @(0.9951252666034236, 0.1652708048635887, 0.1542411818736156, 2.063617003521045, 0.002509536457517263, 0.2360330373669441, 0.5483408417044018, 0.3407535211025703, 5.367398307555874),
    @(0.9637510167226537, 0.1651111235729772, 0.1523547651735839, 2.066133806751765, 0.002511188706227837, 0.2366117307064783, 0.5215104105373598, 0.3363650799180107, 5.380534737484879),
    @(0.9445758383662337, 0.1650165142779443, 0.1512107939044256, 2.067920548614367, 0.002512225233757093, 0.2369887171531389, 0.5050666977532785, 0.3337054845268028, 5.389239590557679),
    @(0.9335734531209994, 0.1649633354491229, 0.1505577173543315, 2.069037499676696, 0.00251282986159746, 0.2372136340173405, 0.4956147311379766, 0.3321877699148956, 5.394483881037075),
    @(0.929852838354293, 0.164945543275735, 0.1503374416081158, 2.069431016031089, 0.002513036030289411, 0.2372911762198981, 0.492415509361507, 0.3316759655205175, 5.396300302147381),
    @(0.9466143140874976, 0.165026457632159, 0.1513320644329355, 2.067721109469673, 0.002512114020222671, 0.2369477497913124, 0.5068165399301563, 0.3339873605652315, 5.388288366130638),
    @(1.003247513277188, 0.1653130615904601, 0.154732305383817, 2.063042227436242, 0.002509116951900409, 0.2358904181183092, 0.5552726942026425, 0.3418965637787039, 5.364206434681137),
    @(1.040492892581824, 0.1655112532606893, 0.1569977423397262, 2.060778150668568, 0.002507232987860705, 0.2352712139666622, 0.5869913272315443, 0.3471717259177751, 5.350581094617127),
    @(1.020593326870596, 0.1654044861020196, 0.1557847089125772, 2.061913944989797, 0.002508231670848939, 0.2355951310582896, 0.5700581329967065, 0.3443466155995338, 5.357659859977957),
    @(0.9456926516260182, 0.1650219584608763, 0.1512772236690978, 2.067810996392353, 0.002512164272663464, 0.2369662456467161, 0.5060254308110927, 0.3338598891313325, 5.388717667851097),
    @(0.8659225809650195, 0.1646567548666695, 0.1466029933380781, 2.077585256915413, 0.002516729322632884, 0.2387543900651039, 0.4371880982680523, 0.3230081004588499, 5.431298910570803)
)

$cols = @("B","C","D","F","G","J","K","L","O")
$startRow = 2

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range($cols[$j] + $row).Value = $data[$i][$j]
    }
}